# Auto-generated edit script applying numeric updates to Kujata_Profits.xlsx
# as described by the commit diff (scheduled runner data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3700.2
$ws.Range("I18").Value = 5500.5
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 5500.5
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = -5216.5
$ws.Range("N18").Value = -3068

$ws.Range("H74").Value = 3650
$ws.Range("I74").Value = 3650
$ws.Range("K74").Value = 3650
$ws.Range("M74").Value = -2714

$ws.Range("H77").Value = 3650
$ws.Range("I77").Value = 3650
$ws.Range("K77").Value = 18250
$ws.Range("M77").Value = -13570

$ws.Range("H132").Value = 6668951
$ws.Range("I132").Value = 9011259
$ws.Range("K132").Value = 27033777
$ws.Range("M132").Value = -27031247

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H138").Value = 1037154.94
$ws.Range("I138").Value = 1761.875
$ws.Range("J138").Value = 1280776.9
$ws.Range("K138").Value = 5285.625
$ws.Range("L138").Value = 3842330.7
$ws.Range("M138").Value = -145.625
$ws.Range("N138").Value = -3852610.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3404.4167
$ws.Range("I32").Value = 3667.6
$ws.Range("J32").Value = 2088.5
$ws.Range("K32").Value = 3667.6
$ws.Range("L32").Value = 2088.5
$ws.Range("M32").Value = -3380.6
$ws.Range("N32").Value = -2662.5

$ws.Range("H74").Value = 1350
$ws.Range("I74").Value = 1597.6666
$ws.Range("J74").Value = 607
$ws.Range("K74").Value = 1597.6666
$ws.Range("L74").Value = 607
$ws.Range("M74").Value = -723.6666
$ws.Range("N74").Value = -2355

$ws.Range("H77").Value = 1350
$ws.Range("I77").Value = 1597.6666
$ws.Range("J77").Value = 607
$ws.Range("K77").Value = 7988.333000000001
$ws.Range("L77").Value = 3035
$ws.Range("M77").Value = -3620.333000000001
$ws.Range("N77").Value = -11771

$ws.Range("H132").Value = 3033.081
$ws.Range("I132").Value = 2765.4827
$ws.Range("K132").Value = 8296.4481
$ws.Range("M132").Value = -5766.4481

$ws.Range("H135").Value = 17203.5
$ws.Range("J135").Value = 17203.5
$ws.Range("L135").Value = 17203.5
$ws.Range("N135").Value = -27343.5

$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4916.8076
$ws.Range("I134").Value = 1031.4783
$ws.Range("K134").Value = 3094.4349
$ws.Range("M134").Value = -559.4349000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1087.4028
$ws.Range("I31").Value = 745.4231
$ws.Range("K31").Value = 745.4231
$ws.Range("M31").Value = -450.4231

$ws.Range("H34").Value = 1087.4028
$ws.Range("I34").Value = 745.4231
$ws.Range("K34").Value = 745.4231
$ws.Range("M34").Value = -543.4231

$ws.Range("H134").Value = 6945663
$ws.Range("I134").Value = 9010105
$ws.Range("K134").Value = 27030315
$ws.Range("M134").Value = -27027780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3033.3333
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 3218.182
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 9654.545999999998
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -9910.545999999998

$ws.Range("H64").Value = 4492.3076
$ws.Range("J64").Value = 4783.3335
$ws.Range("L64").Value = 14350.0005
$ws.Range("N64").Value = -14890.0005

$ws.Range("H67").Value = 4492.3076
$ws.Range("J67").Value = 4783.3335
$ws.Range("L67").Value = 14350.0005
$ws.Range("N67").Value = -16222.0005

$ws.Range("H68").Value = 1852.9344
$ws.Range("I68").Value = 600
$ws.Range("J68").Value = 1873.8167
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 5621.4501
$ws.Range("M68").Value = -989
$ws.Range("N68").Value = -7243.4501

$ws.Range("H71").Value = 1852.9344
$ws.Range("I71").Value = 600
$ws.Range("J71").Value = 1873.8167
$ws.Range("K71").Value = 5400
$ws.Range("L71").Value = 16864.3503
$ws.Range("M71").Value = -1344
$ws.Range("N71").Value = -24976.3503

$ws.Range("H76").Value = 6615.5
$ws.Range("I76").Value = 5056.5
$ws.Range("J76").Value = 6838.2144
$ws.Range("K76").Value = 15169.5
$ws.Range("L76").Value = 20514.6432
$ws.Range("M76").Value = -14786.5
$ws.Range("N76").Value = -21280.6432

$ws.Range("H79").Value = 6615.5
$ws.Range("I79").Value = 5056.5
$ws.Range("J79").Value = 6838.2144
$ws.Range("K79").Value = 15169.5
$ws.Range("L79").Value = 20514.6432
$ws.Range("M79").Value = -13843.5
$ws.Range("N79").Value = -23166.6432

$ws.Range("H88").Value = 3584.375
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 7891.6665
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 23674.9995
$ws.Range("M88").Value = -2572
$ws.Range("N88").Value = -24530.9995

$ws.Range("H91").Value = 3584.375
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 7891.6665
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 23674.9995
$ws.Range("M91").Value = -1518
$ws.Range("N91").Value = -26638.9995

$ws.Range("H94").Value = 4566.6665
$ws.Range("I94").Value = 2750
$ws.Range("J94").Value = 5475
$ws.Range("K94").Value = 8250
$ws.Range("L94").Value = 16425
$ws.Range("M94").Value = -7574
$ws.Range("N94").Value = -17777

$ws.Range("H100").Value = 3150.4443
$ws.Range("J100").Value = 3150.4443
$ws.Range("L100").Value = 9451.332900000001
$ws.Range("N100").Value = -11073.3329

$ws.Range("H106").Value = 4289.857
$ws.Range("J106").Value = 4289.857
$ws.Range("L106").Value = 12869.571
$ws.Range("N106").Value = -14761.571

$ws.Range("H107").Value = 6022.8423
$ws.Range("J107").Value = 7905.7856
$ws.Range("L107").Value = 23717.3568
$ws.Range("N107").Value = -27557.3568

$ws.Range("H109").Value = 93614
$ws.Range("I109").Value = 200830.8
$ws.Range("J109").Value = 4266.6665
$ws.Range("K109").Value = 602492.3999999999
$ws.Range("L109").Value = 12799.9995
$ws.Range("M109").Value = -601452.3999999999
$ws.Range("N109").Value = -14879.9995

$ws.Range("H112").Value = 14988
$ws.Range("J112").Value = 53500
$ws.Range("L112").Value = 160500
$ws.Range("N112").Value = -162716

$ws.Range("H118").Value = 1000
$ws.Range("I118").Value = 1000
$ws.Range("K118").Value = 3000
$ws.Range("M118").Value = -1757

$ws.Range("H122").Value = 581.3
$ws.Range("I122").Value = 427.5
$ws.Range("K122").Value = 3847.5
$ws.Range("M122").Value = -1397.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17311790
$ws.Range("I70").Value = 19234350
$ws.Range("J70").Value = 15389229
$ws.Range("K70").Value = 19234350
$ws.Range("L70").Value = 15389229
$ws.Range("M70").Value = -19234080
$ws.Range("N70").Value = -15389769

$ws.Range("H73").Value = 17311790
$ws.Range("I73").Value = 19234350
$ws.Range("J73").Value = 15389229
$ws.Range("K73").Value = 19234350
$ws.Range("L73").Value = 15389229
$ws.Range("M73").Value = -19233414
$ws.Range("N73").Value = -15391101

$ws.Range("H102").Value = 1486.1818
$ws.Range("I102").Value = 2217.8333
$ws.Range("J102").Value = 1068.0952
$ws.Range("K102").Value = 2217.8333
$ws.Range("L102").Value = 1068.0952
$ws.Range("M102").Value = -595.8332999999998
$ws.Range("N102").Value = -4312.0952

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1437.5
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 1940
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 1940
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -2316

$ws.Range("H100").Value = 2267.6667
$ws.Range("I100").Value = 2001.5
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 2001.5
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -1460.5
$ws.Range("N100").Value = -3882

$ws.Range("H122").Value = 20243030
$ws.Range("I122").Value = 28336344
$ws.Range("K122").Value = 85009032
$ws.Range("M122").Value = -85006582

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 8325
$ws.Range("J86").Value = 8325
$ws.Range("L86").Value = 8325
$ws.Range("N86").Value = -10571

$ws.Range("H89").Value = 8325
$ws.Range("J89").Value = 8325
$ws.Range("L89").Value = 41625
$ws.Range("N89").Value = -52857
